# Add a new hidden worksheet "Лаба1 задания" with a small table,
# and add a note to the existing sheet's B18 cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing sheet (Лист1) ---
# B18 gets a new text value referencing the task note.
$ws1.Range("B18").Value = "Задание 3_8 без функции Math.Pow"

# Update the selected cell on the visible sheet.
$ws1.Range("A16").Select()

# Adjust column widths: column B becomes its own (slightly narrower) bestFit
# width once B18 carries text, while column C keeps its previous width.
$ws1.Columns.Item(2).ColumnWidth = 8

# --- Add the new hidden worksheet ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Лаба1 задания"

$ws2.Range("A1").Value = "ФИО"
$ws2.Range("B1").Value = "Задание"

$ws2.Range("A2").Value = "Родина Ксения Витальевна"
$ws2.Range("B2").Value = "Задание 1_14 и 1_15 без использования кортежей"

$ws2.Range("A3").Value = "Оганезов Михаил Алексеевич"
$ws2.Range("B3").Value = "Задание 1_14 и 1_15 без использования временной переменной"

$ws2.Columns.Item(1).ColumnWidth = 32.333333333333336
$ws2.Columns.Item(2).ColumnWidth = 58.833333333333336

$ws2.Range("B4").Select()

$ws2.Visible = $false

$ws1.Activate()
